# This script updates the "Max Cr" / "Max Date" values for the first four
# wells on both the "Regional for Mapping" sheet (columns V/W) and the
# "Regional Exhibit" sheet (columns G/H), per the recalculated maxima.
#
# Because these replacement values look numeric / date-like, a plain
# Range.Value assignment would make Excel auto-convert them into a real
# number / date serial (changing both the stored type and the cell style,
# since a date needs a date number-format). The source file stores these
# values as plain text, so we explicitly:
#   1. stash the cell's current (pristine) format in a scratch cell,
#   2. force the cell to Text format ("@") before assigning the value so
#      the string is kept verbatim,
#   3. copy the pristine format back from the scratch cell so the cell's
#      style index / appearance is unchanged,
#   4. clear the scratch cell.

$wb = $excel.ActiveWorkbook

function Set-TextValue($Sheet, $CellAddress, $Text, $Scratch) {
    $cell = $Sheet.Range($CellAddress)

    # Preserve the cell's current formatting.
    $cell.Copy()
    $Scratch.PasteSpecial(-4122)  # xlPasteFormats

    # Force text storage so numeric/date-looking strings stay as text.
    $cell.NumberFormat = "@"
    $cell.Value = $Text

    # Restore the original formatting/style.
    $Scratch.Copy()
    $cell.PasteSpecial(-4122)  # xlPasteFormats

    $Scratch.Clear()
}

# "Regional for Mapping" sheet - Max Cr (V) and Max Date (W) columns
$wsMap = $wb.Worksheets.Item("Regional for Mapping")
$scratchMap = $wsMap.Range("ZZ1")

Set-TextValue $wsMap "V9"  "2.4"        $scratchMap
Set-TextValue $wsMap "W9"  "2005-03-23" $scratchMap

Set-TextValue $wsMap "V10" "1.4"        $scratchMap
Set-TextValue $wsMap "W10" "2005-03-22" $scratchMap

Set-TextValue $wsMap "V11" "3.6"        $scratchMap
Set-TextValue $wsMap "W11" "2005-03-21" $scratchMap

Set-TextValue $wsMap "V12" "2.74"       $scratchMap
Set-TextValue $wsMap "W12" "2002-05-17" $scratchMap

# "Regional Exhibit" sheet - Max Cr (G) and Max Date (H) columns
$wsExhibit = $wb.Worksheets.Item("Regional Exhibit")
$scratchExhibit = $wsExhibit.Range("ZZ1")

Set-TextValue $wsExhibit "G13" "2.4"     $scratchExhibit
Set-TextValue $wsExhibit "H13" "3/23/05" $scratchExhibit

Set-TextValue $wsExhibit "G14" "1.4"     $scratchExhibit
Set-TextValue $wsExhibit "H14" "3/22/05" $scratchExhibit

Set-TextValue $wsExhibit "G15" "3.6"     $scratchExhibit
Set-TextValue $wsExhibit "H15" "3/21/05" $scratchExhibit

Set-TextValue $wsExhibit "G16" "2.74"    $scratchExhibit
Set-TextValue $wsExhibit "H16" "5/17/02" $scratchExhibit
